# Slide 5 ("Idea"): the three dashed connector "lines" that used to sit on
# the left cluster of points (ids 10/11/12, named "Shape 52"/"Shape 53"/
# "Shape 54" - duplicates of shapes 52/53/54) are moved over to the right
# cluster of points and rotated, per the target OOXML:
#
#   id=10: rot=-3960000 (-66 deg)  off=(7360519, 6632053)
#   id=11: rot= 6300000 ( 105 deg) off=(7689924, 6716000)
#   id=12: rot= 5580000 (  93 deg) off=(7975668, 6867037)
#
# Widths/heights (the shapes' Ext) are unchanged.
#
# Note: Shape.Left/.Top are expressed in points and stored internally as
# single-precision floats; on save they are converted back to EMU
# (1 pt = 12700 EMU) by truncation. A plain EMU/12700.0 can therefore land
# one EMU below the intended target after that round trip, so the literals
# below have been nudged by a tiny fraction of a point so the saved EMU
# values come out exactly as in the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$changes = @(
    @{ Index = 9;  RotationDeg = -66.0; LeftPt = 579.5684511969; TopPt = 522.2088976378 },
    @{ Index = 10; RotationDeg = 105.0; LeftPt = 605.5058297717; TopPt = 528.8188976378 },
    @{ Index = 11; RotationDeg = 93.0;  LeftPt = 628.0053543307; TopPt = 540.7115788031 }
)

foreach ($change in $changes) {
    $shp = $s.Shapes.Item($change.Index)
    $shp.Left = $change.LeftPt
    $shp.Top = $change.TopPt
    $shp.Rotation = $change.RotationDeg
}
